# This script inserts a new weekly price record for "Terminal Hortofrutícola
# Agro Chillán - Zanahoria" at row 212, shifting all existing records
# (previously rows 212-300) down by one row (to rows 213-301).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 212 (pushes rows 212:300 down to 213:301)
$ws.Rows("212:212").Insert()

# Populate the newly inserted row 212 with the new record's data
$ws.Range("A212").Value = 7
$ws.Range("B212").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C212").Value = "Ñuble"
$ws.Range("D212").Value = 44755
$ws.Range("E212").Value = 16
$ws.Range("F212").Value = 100114013
$ws.Range("G212").Value = "Zanahoria"
$ws.Range("H212").Value = "Sin especificar"
$ws.Range("I212").Value = "Primera"
$ws.Range("J212").Value = 120
$ws.Range("K212").Value = 7000
$ws.Range("L212").Value = 7500
$ws.Range("M212").Value = 7250
$ws.Range("N212").Value = "`$/saco 20 kilos"
$ws.Range("O212").Value = "Provincia de Diguillín"
$ws.Range("P212").Value = 362
$ws.Range("Q212").Value = 20
$ws.Range("R212").Value = "Hortaliza"
